$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6799.875
$ws.Range("I18").Value = 624.75
$ws.Range("K18").Value = 624.75
$ws.Range("M18").Value = -340.75

$ws.Range("H40").Value = 4030.6875
$ws.Range("J40").Value = 4458
$ws.Range("L40").Value = 4458
$ws.Range("N40").Value = -4808

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H87").Value = 32920.715

$ws.Range("H90").Value = 32920.715

$ws.Range("H111").Value = 1656.4
$ws.Range("J111").Value = 1444.1666
$ws.Range("L111").Value = 4332.4998
$ws.Range("N111").Value = -10466.4998

$ws.Range("H138").Value = 2038.5758
$ws.Range("I138").Value = 1237.7142
$ws.Range("J138").Value = 3440.0833
$ws.Range("K138").Value = 3713.1426
$ws.Range("L138").Value = 10320.2499
$ws.Range("M138").Value = 1426.8574
$ws.Range("N138").Value = -20600.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 166.5
$ws.Range("I5").Value = 149.75
$ws.Range("K5").Value = 149.75
$ws.Range("M5").Value = -37.75

$ws.Range("H21").Value = 599.75
$ws.Range("I21").Value = 149
$ws.Range("J21").Value = 750
$ws.Range("K21").Value = 149
$ws.Range("L21").Value = 750
$ws.Range("M21").Value = 225
$ws.Range("N21").Value = -1498

$ws.Range("H32").Value = 3758749
$ws.Range("I32").Value = 718848.9399999999
$ws.Range("J32").Value = 22335916
$ws.Range("K32").Value = 718848.9399999999
$ws.Range("L32").Value = 22335916
$ws.Range("M32").Value = -718561.9399999999
$ws.Range("N32").Value = -22336490

$ws.Range("H37").Value = 14647.5
$ws.Range("J37").Value = 25625
$ws.Range("L37").Value = 25625
$ws.Range("N37").Value = -26171

$ws.Range("H110").Value = 1243.55
$ws.Range("I110").Value = 993.2632
$ws.Range("K110").Value = 993.2632
$ws.Range("M110").Value = 1051.7368

$ws.Range("H122").Value = 2417.652
$ws.Range("I122").Value = 1701.3889
$ws.Range("J122").Value = 4996.2
$ws.Range("K122").Value = 5104.1667
$ws.Range("L122").Value = 14988.6
$ws.Range("M122").Value = -2654.1667
$ws.Range("N122").Value = -19888.6

$ws.Range("H124").Value = 63429
$ws.Range("J124").Value = 63429
$ws.Range("L124").Value = 63429
$ws.Range("N124").Value = -73249

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 166.5
$ws.Range("I4").Value = 149.75
$ws.Range("K4").Value = 149.75
$ws.Range("M4").Value = -34.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1932.1578
$ws.Range("I16").Value = 1573.4667
$ws.Range("K16").Value = 1573.4667
$ws.Range("M16").Value = -1286.4667

$ws.Range("H22").Value = 1328.4
$ws.Range("I22").Value = 1160.5
$ws.Range("K22").Value = 1160.5
$ws.Range("M22").Value = -810.5

$ws.Range("H58").Value = 889
$ws.Range("I58").Value = 843
$ws.Range("K58").Value = 843
$ws.Range("M58").Value = -640

$ws.Range("H113").Value = 1932.1578
$ws.Range("I113").Value = 1573.4667
$ws.Range("K113").Value = 1573.4667
$ws.Range("M113").Value = 596.5333000000001

$ws.Range("H118").Value = 99999.5
$ws.Range("J118").Value = 99999.5
$ws.Range("L118").Value = 99999.5
$ws.Range("N118").Value = -103313.5

$ws.Range("H136").Value = 889
$ws.Range("I136").Value = 843
$ws.Range("K136").Value = 2529
$ws.Range("M136").Value = 21

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3977
$ws.Range("I3").Value = 1015
$ws.Range("J3").Value = 9901
$ws.Range("K3").Value = 3045
$ws.Range("L3").Value = 29703
$ws.Range("M3").Value = -2933
$ws.Range("N3").Value = -29927

$ws.Range("H118").Value = 115560.89
$ws.Range("I118").Value = 166841.17
$ws.Range("K118").Value = 500523.51
$ws.Range("M118").Value = -499280.51

$ws.Range("H129").Value = 225475.22
$ws.Range("I129").Value = 500304.75
$ws.Range("K129").Value = 1500914.25
$ws.Range("M129").Value = -1495914.25

$ws.Range("H131").Value = 17941996
$ws.Range("I131").Value = 41750960
$ws.Range("K131").Value = 125252880
$ws.Range("M131").Value = -125247840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H7").Value = 7909.579
$ws.Range("I7").Value = 5129.077
$ws.Range("K7").Value = 5129.077
$ws.Range("M7").Value = -5017.077

$ws.Range("H46").Value = 23447.736
$ws.Range("I46").Value = 31279.072
$ws.Range("J46").Value = 1520
$ws.Range("K46").Value = 31279.072
$ws.Range("L46").Value = 1520
$ws.Range("M46").Value = -31091.072
$ws.Range("N46").Value = -1896

$ws.Range("H61").Value = 58185.848
$ws.Range("I61").Value = 53595.55
$ws.Range("J61").Value = 73486.836
$ws.Range("K61").Value = 53595.55
$ws.Range("L61").Value = 73486.836
$ws.Range("M61").Value = -53393.55
$ws.Range("N61").Value = -73890.836

$ws.Range("H93").Value = 20778.588
$ws.Range("I93").Value = 1189.5
$ws.Range("J93").Value = 112194.336
$ws.Range("K93").Value = 1189.5
$ws.Range("L93").Value = 112194.336
$ws.Range("M93").Value = 58.5
$ws.Range("N93").Value = -114690.336

$ws.Range("H113").Value = 58185.848
$ws.Range("I113").Value = 53595.55
$ws.Range("J113").Value = 73486.836
$ws.Range("K113").Value = 53595.55
$ws.Range("L113").Value = 73486.836
$ws.Range("M113").Value = -51425.55
$ws.Range("N113").Value = -77826.836

$ws.Range("H126").Value = 7909.579
$ws.Range("I126").Value = 5129.077
$ws.Range("K126").Value = 15387.231
$ws.Range("M126").Value = -12917.231

$ws.Range("H132").Value = 3212.1462
$ws.Range("I132").Value = 2787.9285
$ws.Range("J132").Value = 4125.846
$ws.Range("K132").Value = 8363.7855
$ws.Range("L132").Value = 12377.538
$ws.Range("M132").Value = -5833.7855
$ws.Range("N132").Value = -17437.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 21765158
$ws.Range("I107").Value = 1447.8422
$ws.Range("K107").Value = 4343.5266
$ws.Range("M107").Value = -2423.5266

$ws.Range("H136").Value = 4802.146
$ws.Range("I136").Value = 3337.465
$ws.Range("K136").Value = 10012.395
$ws.Range("M136").Value = -7462.395
